$d = $word.ActiveDocument
$sec = $d.Sections(1)

# The header/footer logo pictures were saved with a mismatched display
# name: the Pearson logo (footers, primary + first page) is labelled
# "image1.png" and the BTEC logo (headers, primary + first page) is
# labelled "image2.jpg" - the reverse of what they should be. Rename
# each InlineShape so the name matches its actual media part
# (footers -> image2.png, headers -> image1.jpg).

for ($i = 1; $i -le 2; $i++) {
    $footerShape = $sec.Footers($i).Range.InlineShapes(1)
    $footerShape.Name = "image2.png"

    $headerShape = $sec.Headers($i).Range.InlineShapes(1)
    $headerShape.Name = "image1.jpg"
}
